$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - previously only held the subscription user-story text in A11 with
# the rest of the row blank; now the whole row is filled out with Joel's
# rewritten user story plus the new "Must/Should/Could/Won't Have" items.
# B11:E11 pick up the same vertical-centered / wrap-text formatting already
# used by the rest of this table block (rows 12 & 14).
$ws.Range("B11:E11").VerticalAlignment = -4108
$ws.Range("B11:E11").WrapText = $true
$ws.Range("A11").Value = "I, Joel, am a Twitch user called MasterTwitchUser that has just subscribed to the channel."
$ws.Range("B11").Value = "Bot doesn't recognize a user trying to imitate a sub using a chat message."
$ws.Range("C11").Value = "Should have the ablity to work on a Twitch Prime sub and other sub types."
$ws.Range("D11").Value = "A bot message congratulating the user for subscribing."
$ws.Range("E11").Value = "Give the user coins of points for subscribing."

# Row 12 - the "Must Have" cell text stays the same; the others get typo
# fixes, and the trailing "Joel" note in F12 is removed.
$ws.Range("C12").Value = "Should log a message out of the chats view."
$ws.Range("D12").Value = "Logs name and can stores the subscriptions in a data file for mining purposes."
$ws.Range("E12").Value = "Won't print out a message to the views showing the user subscribed because that is already built into Twitch."
$ws.Range("F12").ClearContents()

# Row 13 - previously only A13 held the gambling user-story text; now the
# whole row is filled out with Joel's rewritten story plus new items.
# B13:E13 pick up the same vertical-centered / wrap-text formatting already
# used by the rest of this table block (rows 12 & 14).
$ws.Range("B13:E13").VerticalAlignment = -4108
$ws.Range("B13:E13").WrapText = $true
$ws.Range("A13").Value = "I, Joel, am a Twitch user TheGreatGabby01 who wants to gamble some of their coins for potentailly more coins. I can type ""!gamble xxx"" in chat. xxx = heads/tails"
$ws.Range("B13").Value = "Stop the user from gambling if they don't have enough coins."
$ws.Range("C13").Value = "The ability to check, add, and deduct points from the user."
$ws.Range("D13").Value = "Other people placing bets on that person's bet."
$ws.Range("E13").Value = "Won't have any monetary value because I believe that's illegal."

# Row 14 - B14:E14 content is unchanged; only the trailing "Joel" note in
# F14 is removed.
$ws.Range("F14").ClearContents()

# Update the view so the saved selection matches the author's editing spot.
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 11
